$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 463, pushing existing rows 463-487 down to 464-488.
# Excel's row-insert copies formatting from the row above, which already
# gives the date cell (column D) the "YYYY-MM-DD HH:MM:SS" custom style
# used throughout the column.
$ws.Rows.Item(463).Insert()

# Populate the newly-inserted row 463 with the new weekly record.
$ws.Cells.Item(463, 1).Value = 8
$ws.Cells.Item(463, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(463, 3).Value = "Coquimbo"
$ws.Cells.Item(463, 4).Value = 45041
$ws.Cells.Item(463, 5).Value = 4
$ws.Cells.Item(463, 6).Value = 100112032
$ws.Cells.Item(463, 7).Value = "Zapallo italiano"
$ws.Cells.Item(463, 8).Value = "Sin especificar"
$ws.Cells.Item(463, 9).Value = "Primera"
$ws.Cells.Item(463, 10).Value = 400
$ws.Cells.Item(463, 11).Value = 9000
$ws.Cells.Item(463, 12).Value = 10000
$ws.Cells.Item(463, 13).Value = 9500
$ws.Cells.Item(463, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(463, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(463, 16).Value = 136
$ws.Cells.Item(463, 17).Value = 70
$ws.Cells.Item(463, 18).Value = "Hortaliza"
